# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for Alcachofa (Comercializadora del Agro de
# Limarí) above the existing block of rows 190-195, pushing that block down
# to rows 193-198 (dimension grows from A1:R195 to A1:R198).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 190:195 down to 193:198 by inserting 3 blank rows.
$ws.Rows("190:192").Insert()

# Helper values shared by every row in this subset.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112013
$categoria = "Alcachofa"
$clasif    = "Hortaliza"

# --- New row 190: Española / Primera -----------------------------------
$ws.Cells.Item(190, 1).Value  = $mercadoId
$ws.Cells.Item(190, 2).Value  = $mercado
$ws.Cells.Item(190, 3).Value  = $region
$ws.Cells.Item(190, 4).Value  = 44714
$ws.Cells.Item(190, 5).Value  = $codreg
$ws.Cells.Item(190, 6).Value  = $catId
$ws.Cells.Item(190, 7).Value  = $categoria
$ws.Cells.Item(190, 8).Value  = "Española"
$ws.Cells.Item(190, 9).Value  = "Primera"
$ws.Cells.Item(190, 10).Value = 500
$ws.Cells.Item(190, 11).Value = 16000
$ws.Cells.Item(190, 12).Value = 17000
$ws.Cells.Item(190, 13).Value = 16500
$ws.Cells.Item(190, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(190, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(190, 16).Value = 550
$ws.Cells.Item(190, 17).Value = 30
$ws.Cells.Item(190, 18).Value = $clasif

# --- New row 191: Española / Segunda ------------------------------------
$ws.Cells.Item(191, 1).Value  = $mercadoId
$ws.Cells.Item(191, 2).Value  = $mercado
$ws.Cells.Item(191, 3).Value  = $region
$ws.Cells.Item(191, 4).Value  = 44714
$ws.Cells.Item(191, 5).Value  = $codreg
$ws.Cells.Item(191, 6).Value  = $catId
$ws.Cells.Item(191, 7).Value  = $categoria
$ws.Cells.Item(191, 8).Value  = "Española"
$ws.Cells.Item(191, 9).Value  = "Segunda"
$ws.Cells.Item(191, 10).Value = 360
$ws.Cells.Item(191, 11).Value = 14000
$ws.Cells.Item(191, 12).Value = 15000
$ws.Cells.Item(191, 13).Value = 14500
$ws.Cells.Item(191, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(191, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(191, 16).Value = 362
$ws.Cells.Item(191, 17).Value = 40
$ws.Cells.Item(191, 18).Value = $clasif

# --- New row 192: Madrigal / Primera ------------------------------------
$ws.Cells.Item(192, 1).Value  = $mercadoId
$ws.Cells.Item(192, 2).Value  = $mercado
$ws.Cells.Item(192, 3).Value  = $region
$ws.Cells.Item(192, 4).Value  = 44714
$ws.Cells.Item(192, 5).Value  = $codreg
$ws.Cells.Item(192, 6).Value  = $catId
$ws.Cells.Item(192, 7).Value  = $categoria
$ws.Cells.Item(192, 8).Value  = "Madrigal"
$ws.Cells.Item(192, 9).Value  = "Primera"
$ws.Cells.Item(192, 10).Value = 140
$ws.Cells.Item(192, 11).Value = 16000
$ws.Cells.Item(192, 12).Value = 17000
$ws.Cells.Item(192, 13).Value = 16500
$ws.Cells.Item(192, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(192, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(192, 16).Value = 412
$ws.Cells.Item(192, 17).Value = 40
$ws.Cells.Item(192, 18).Value = $clasif
